# Update "想去人数" (want-to-go count) figures in F column across sheets
# 展览 (Exhibition), 演出 (Performance), 全部类型 (All types).
# 本地生活 (Local life) sheet is unchanged.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 26     # was 25
$ws.Range("F6").Value = 3001   # was 3000
$ws.Range("F7").Value = 1282   # was 1281
$ws.Range("F8").Value = 445    # was 443
$ws.Range("F11").Value = 20    # was 17
$ws.Range("F12").Value = 742   # was 741
$ws.Range("F13").Value = 294   # was 292
$ws.Range("F18").Value = 124   # was 123
$ws.Range("F22").Value = 6768  # was 6763
$ws.Range("F25").Value = 465   # was 464
$ws.Range("F26").Value = 1257  # was 1255
$ws.Range("F27").Value = 6260  # was 6259
$ws.Range("F30").Value = 1866  # was 1867
$ws.Range("F31").Value = 6013  # was 6010
$ws.Range("F36").Value = 428   # was 427
$ws.Range("F37").Value = 4610  # was 4524
$ws.Range("F39").Value = 190   # was 189
$ws.Range("F43").Value = 2414  # was 2413
$ws.Range("F48").Value = 347   # was 343
$ws.Range("F49").Value = 2070  # was 2067
$ws.Range("F50").Value = 10    # was 8

# --- 演出 sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 196    # was 195
$ws.Range("F7").Value = 30     # was 29

# --- 全部类型 sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 26     # was 25
$ws.Range("F5").Value = 3001   # was 3000
$ws.Range("F6").Value = 1282   # was 1281
$ws.Range("F7").Value = 445    # was 443
$ws.Range("F10").Value = 196   # was 195
$ws.Range("F12").Value = 294   # was 292
$ws.Range("F17").Value = 124   # was 123
$ws.Range("F21").Value = 6768  # was 6763
$ws.Range("F24").Value = 465   # was 464
$ws.Range("F25").Value = 1257  # was 1255
$ws.Range("F27").Value = 6260  # was 6259
$ws.Range("F29").Value = 1866  # was 1867
$ws.Range("F31").Value = 6013  # was 6010
$ws.Range("F37").Value = 428   # was 427
$ws.Range("F38").Value = 4610  # was 4525
$ws.Range("F40").Value = 190   # was 189
$ws.Range("F44").Value = 2414  # was 2413
$ws.Range("F49").Value = 347   # was 343

$wb.Save()
